$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 173, shifting existing rows 173:248 down to 174:249
$ws.Rows.Item(173).Insert()

# Populate the new row 173 with the new record's data
$ws.Range("A173").Value = 5
$ws.Range("B173").Value = "Macroferia Regional de Talca"
$ws.Range("C173").Value = "Maule"
$ws.Range("D173").Value = 44704
$ws.Range("E173").Value = 7
$ws.Range("F173").Value = 100112009
$ws.Range("G173").Value = "Acelga"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 500
$ws.Range("K173").Value = 3000
$ws.Range("L173").Value = 3000
$ws.Range("M173").Value = 3000
$ws.Range("N173").Value = "$/docena de atados (4 kilos)"
$ws.Range("O173").Value = "Región del Maule"
$ws.Range("P173").Value = 750
$ws.Range("Q173").Value = 4
$ws.Range("R173").Value = "Hortaliza"
